# Insert a new record at row 482 (weekly Perejil price data), pushing the
# existing rows 482-526 down to 483-527.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(482).Insert()

$ws.Cells.Item(482, 1).Value = 6
$ws.Cells.Item(482, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(482, 3).Value = "Metropolitana"
$ws.Cells.Item(482, 4).Value = "2022-07-04"
$ws.Cells.Item(482, 5).Value = 13
$ws.Cells.Item(482, 6).Value = 100112044
$ws.Cells.Item(482, 7).Value = "Perejil"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Primera"
$ws.Cells.Item(482, 10).Value = 120
$ws.Cells.Item(482, 11).Value = 17000
$ws.Cells.Item(482, 12).Value = 18000
$ws.Cells.Item(482, 13).Value = 17417
$ws.Cells.Item(482, 14).Value = "`$/docena de atados"
$ws.Cells.Item(482, 15).Value = "Región Metropolitana"
$ws.Cells.Item(482, 16).Value = 5806
$ws.Cells.Item(482, 17).Value = 3
$ws.Cells.Item(482, 18).Value = "Hortaliza"
